$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.795.54'
$ws.Range('E2').Value = '  -0.55%  '
$ws.Range('D3').Value = '1.933.47'
$ws.Range('E3').Value = '  -0.89%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '242.68'
$ws.Range('E5').Value = '  -1.00%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4888'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2948'
$ws.Range('E8').Value = '  -0.58%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.06863'
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.22'
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '104.70'
$ws.Range('E11').Value = '  -2.50%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07789'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.930.50'
$ws.Range('E13').Value = '  -1.10%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.321'
$ws.Range('E14').Value = '  -2.39%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.7005'
$ws.Range('E15').Value = '  -0.50%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '271.92'
$ws.Range('E16').Value = '  -3.86%  '
$ws.Range('D17').Value = '30.786.68'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.000007702'
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.624'
$ws.Range('E19').Value = '  +2.67%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '13.05'
$ws.Range('E20').Value = '  -1.30%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.002'
$ws.Range('E21').Value = '  +0.16%  '
$ws.Range('E22').Value = '  +0.14%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.519'
$ws.Range('E23').Value = '  +0.69%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '9.793'
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '164.66'
$ws.Range('E25').Value = '  -3.02%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '19.50'
$ws.Range('E26').Value = '  -2.29%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.157'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1037'
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.388'
$ws.Range('E29').Value = '  -1.54%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '4.574'
$ws.Range('E30').Value = '  -0.37%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.549'
$ws.Range('E31').Value = '  -1.87%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.390'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04891'
$ws.Range('E33').Value = '  -1.20%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.7532'
$ws.Range('E34').Value = '  -1.63%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.145'
$ws.Range('E35').Value = '  -2.02%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.002'
$ws.Range('E36').Value = '  +0.20%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.735'
$ws.Range('E37').Value = '  +0.34%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.02000'
$ws.Range('E38').Value = '  -0.54%  '
$ws.Range('B39').Value = 'Aave'
$ws.Range('C39').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '79.22'
$ws.Range('E39').Value = '  +7.48%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.661'
$ws.Range('E40').Value = '  -1.44%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '6.452'
$ws.Range('E41').Value = '  -1.86%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.058'
$ws.Range('E42').Value = '  -3.85%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.8930'
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.4433'
$ws.Range('E44').Value = '  -1.04%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '107.95'
$ws.Range('E45').Value = '  -1.53%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '7.936'
$ws.Range('E46').Value = '  -2.17%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.23%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '981.87'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.1241'
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '36.14'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '9.134'
$ws.Range('E51').Value = '  -1.91%  '
